$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new Mac-Address rows (31 and 32), following the same pattern as
# the preceding data rows.
$newRows = @(
    @(10001, 110030, 10030, "eng", $true, "superadmin", "now()", "now()"),
    @(10001, 110031, 10031, "eng", $true, "superadmin", "now()", "now()")
)

$r = 31
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Match the final view state captured in the saved workbook: scrolled down
# with F30 selected.
$ws.Application.Goto($ws.Range("A25"))
$ws.Range("F30").Select()
